$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format before writing, so that numeric-looking
# strings (e.g. "0.9999", "243.38") are preserved as text instead of being
# auto-converted to numbers by Excel, matching the original inlineStr/text cells.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "26.413.91"
$ws.Range("E2").Value = "  -0.39%  "

# Row 3
$ws.Range("D3").Value = "1.723.43"
$ws.Range("E3").Value = "  -0.27%  "

# Row 4
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "243.38"
$ws.Range("E5").Value = "  -0.74%  "

# Row 6
$ws.Range("E6").Value = "  +0.03%  "

# Row 7
$ws.Range("D7").Value = "0.4910"
$ws.Range("E7").Value = "  +2.08%  "

# Row 8
$ws.Range("D8").Value = "0.2614"
$ws.Range("E8").Value = "  -1.94%  "

# Row 9
$ws.Range("D9").Value = "0.06198"
$ws.Range("E9").Value = "  +0.06%  "

# Row 10
$ws.Range("D10").Value = "1.721.12"
$ws.Range("E10").Value = "  -0.41%  "

# Row 11
$ws.Range("D11").Value = "0.07004"
$ws.Range("E11").Value = "  -2.47%  "

# Row 12
$ws.Range("D12").Value = "15.52"
$ws.Range("E12").Value = "  -0.39%  "

# Row 13
$ws.Range("D13").Value = "4.573"
$ws.Range("E13").Value = "  +1.01%  "

# Row 14
$ws.Range("D14").Value = "0.5995"
$ws.Range("E14").Value = "  -1.67%  "

# Row 15
$ws.Range("D15").Value = "77.23"
$ws.Range("E15").Value = "  +0.10%  "

# Row 16
$ws.Range("E16").Value = "  +0.02%  "

# Row 17
$ws.Range("D17").Value = "26.409.21"
$ws.Range("E17").Value = "  -0.40%  "

# Row 18
$ws.Range("D18").Value = "0.9999"
$ws.Range("E18").Value = "  +0.00%  "

# Row 19
$ws.Range("D19").Value = "0.000007148"
$ws.Range("E19").Value = "  +3.03%  "

# Row 20
$ws.Range("D20").Value = "11.36"
$ws.Range("E20").Value = "  -1.39%  "

# Row 21
$ws.Range("D21").Value = "1.944.19"
$ws.Range("E21").Value = "  -0.55%  "

# Row 22
$ws.Range("D22").Value = "4.477"
$ws.Range("E22").Value = "  -1.02%  "

# Row 23
$ws.Range("D23").Value = "8.587"
$ws.Range("E23").Value = "  -2.49%  "

# Row 24
$ws.Range("D24").Value = "5.149"
$ws.Range("E24").Value = "  -1.98%  "

# Row 25
$ws.Range("E25").Value = "  +0.31%  "

# Row 26
$ws.Range("D26").Value = "15.22"
$ws.Range("E26").Value = "  -0.76%  "

# Row 27
$ws.Range("D27").Value = "1.396"
$ws.Range("E27").Value = "  -0.73%  "

# Row 28
$ws.Range("D28").Value = "107.14"
$ws.Range("E28").Value = "  -0.21%  "

# Row 29
$ws.Range("D29").Value = "1.703"
$ws.Range("E29").Value = "  -3.96%  "

# Row 30
$ws.Range("D30").Value = "3.945"
$ws.Range("E30").Value = "  -0.69%  "

# Row 31
$ws.Range("D31").Value = "0.07965"
$ws.Range("E31").Value = "  -0.87%  "

# Row 32
$ws.Range("D32").Value = "3.673"
$ws.Range("E32").Value = "  -0.39%  "

# Row 33
$ws.Range("D33").Value = "0.04548"
$ws.Range("E33").Value = "  +0.78%  "

# Row 34
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").Value = "0.9993"
$ws.Range("E34").Value = "  -0.02%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.600"
$ws.Range("E35").Value = "  -0.66%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "0.9950"
$ws.Range("E36").Value = "  -0.44%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.6259"
$ws.Range("E37").Value = "  +0.02%  "

# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "0.9244"
$ws.Range("E38").Value = "  +1.63%  "

# Row 39
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.393"
$ws.Range("E39").Value = "  +0.78%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "1.948"
$ws.Range("E40").Value = "  -6.08%  "

# Row 41
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "0.9998"
$ws.Range("E41").Value = "  -0.13%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.01484"
$ws.Range("E42").Value = "  -1.30%  "

# Row 43
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "99.90"
$ws.Range("E43").Value = "  -2.41%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.343"
$ws.Range("E44").Value = "  -4.10%  "

# Row 45
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.3845"
$ws.Range("E45").Value = "  -0.79%  "

# Row 46
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "6.722"
$ws.Range("E46").Value = "  -3.53%  "

# Row 47
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.1161"
$ws.Range("E47").Value = "  -1.79%  "

# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05362"
$ws.Range("E48").Value = "  -0.03%  "

# Row 49
$ws.Range("E49").Value = "  -1.18%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.700"
$ws.Range("E50").Value = "  -0.82%  "

# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.241"
$ws.Range("E51").Value = "  -1.03%  "

# Restore column D to the default (unstyled) cell style now that the values
# have been written as text, so no stray style index is left on the cells.
$dRange.Style = "Normal"
